# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    2  = 8356
    3  = 7813
    8  = 128
    9  = 120
    10 = 167
    12 = 710
    13 = 129
    14 = 1350
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
